$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 Col 1: 10 x 91 -> 64 x 96
$range = $t.Cell(1, 1).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>64 x 96</w:t><w:br/><w:t xml:space=`"preserve`">  9    6</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 1 Col 2: 91 x 86 -> 17 x 77
$range = $t.Cell(1, 2).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>17 x 77</w:t><w:br/><w:t xml:space=`"preserve`">  7    7</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 1 Col 3: 41 x 17 -> 73 x 45
$range = $t.Cell(1, 3).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>73 x 45</w:t><w:br/><w:t xml:space=`"preserve`">  4    5</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 2 Col 1: 35 x 85 -> 85 x 75
$range = $t.Cell(2, 1).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>85 x 75</w:t><w:br/><w:t xml:space=`"preserve`">  7    5</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 2 Col 2: 31 x 45 -> 25 x 67
$range = $t.Cell(2, 2).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>25 x 67</w:t><w:br/><w:t xml:space=`"preserve`">  6    7</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 2 Col 3: 39 x 56 -> 75 x 78
$range = $t.Cell(2, 3).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>75 x 78</w:t><w:br/><w:t xml:space=`"preserve`">  7    8</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 3 Col 1: 87 x 19 -> 68 x 78
$range = $t.Cell(3, 1).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>68 x 78</w:t><w:br/><w:t xml:space=`"preserve`">  7    8</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 3 Col 2: 84 x 35 -> 53 x 27
$range = $t.Cell(3, 2).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>53 x 27</w:t><w:br/><w:t xml:space=`"preserve`">  2    7</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 3 Col 3: 63 x 14 -> 20 x 80
$range = $t.Cell(3, 3).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>20 x 80</w:t><w:br/><w:t xml:space=`"preserve`">  8    0</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 4 Col 1: 20 x 57 -> 87 x 70
$range = $t.Cell(4, 1).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>87 x 70</w:t><w:br/><w:t xml:space=`"preserve`">  7    0</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>7|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 4 Col 2: 81 x 65 -> 82 x 65
$range = $t.Cell(4, 2).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>82 x 65</w:t><w:br/><w:t xml:space=`"preserve`">  6    5</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 4 Col 3: 41 x 21 -> 94 x 84
$range = $t.Cell(4, 3).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>94 x 84</w:t><w:br/><w:t xml:space=`"preserve`">  8    4</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 5 Col 1: 51 x 70 -> 64 x 42
$range = $t.Cell(5, 1).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>64 x 42</w:t><w:br/><w:t xml:space=`"preserve`">  4    2</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>6|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 5 Col 2: 39 x 31 -> 54 x 37
$range = $t.Cell(5, 2).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>54 x 37</w:t><w:br/><w:t xml:space=`"preserve`">  3    7</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

# Row 5 Col 3: 10 x 61 -> 96 x 73
$range = $t.Cell(5, 3).Range
$range.InsertXML("<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr><w:t>96 x 73</w:t><w:br/><w:t xml:space=`"preserve`">  7    3</w:t><w:br/><w:t xml:space=`"preserve`">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")
